# Generate Report for Handoff
# Updates the localization-status report: the "b.md" row has been
# re-handed-off. Status moves from "Handed back: in sync with en-US" to
# "Ready for handoff" on the Overview sheet (both locales) and on each
# locale sheet; the locale sheets also get a fresh handoff file name /
# datetime, the Content Duplicate flag flips to False, and a staleness
# warning is recorded in the Error Detail column because the handback
# file has not yet caught up with the new handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 14:41:42"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-22 14:41:37"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1891fc2d1847f0835c25cf17fd7b7a21e16bd720/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b710459ee22751783d683440a65c2c203ffb2fd/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-22 14:41:42"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1891fc2d1847f0835c25cf17fd7b7a21e16bd720/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b710459ee22751783d683440a65c2c203ffb2fd/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
